$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.009.37'
$ws.Range("E2").Value = '  +2.88%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.803.81'
$ws.Range("E3").Value = '  +0.99%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '706.13'
$ws.Range("E5").Value = '  +11.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.02'
$ws.Range("E6").Value = '  +4.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.801.71'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("E10").Value = '  +3.45%  '
$ws.Range("E11").Value = '  +8.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.464'
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000258'
$ws.Range("E13").Value = '  +7.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.47'
$ws.Range("E14").Value = '  +4.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.441.74'
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.799.69'
$ws.Range("E16").Value = '  +0.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.013.89'
$ws.Range("E17").Value = '  +2.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.02'
$ws.Range("E18").Value = '  +2.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.26'
$ws.Range("E19").Value = '  +3.50%  '
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.29'
$ws.Range("E21").Value = '  +19.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '484.56'
$ws.Range("E22").Value = '  +4.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.720'
$ws.Range("E23").Value = '  +2.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.96'
$ws.Range("E24").Value = '  +2.29%  '
$ws.Range("E25").Value = '  +1.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.56'
$ws.Range("E26").Value = '  +3.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.69'
$ws.Range("E27").Value = '  +5.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("E28").Value = '  +3.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.951.93'
$ws.Range("E29").Value = '  +0.83%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.14'
$ws.Range("E31").Value = '  +17.73%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.68'
$ws.Range("E32").Value = '  +8.72%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.31'
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.70'
$ws.Range("E34").Value = '  +4.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.180'
$ws.Range("E35").Value = '  +3.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.33'
$ws.Range("E36").Value = '  +4.88%  '
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.751.67'
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.104'
$ws.Range("E39").Value = '  +3.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.52'
$ws.Range("E40").Value = '  +6.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.04'
$ws.Range("E41").Value = '  +4.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.26'
$ws.Range("E42").Value = '  +14.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000333'
$ws.Range("E43").Value = '  +26.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.972'
$ws.Range("E44").Value = '  +1.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.81'
$ws.Range("E47").Value = '  +7.14%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '160.82'
$ws.Range("E48").Value = '  +2.56%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '49.27'
$ws.Range("E49").Value = '  +4.79%  '
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.301'
$ws.Range("E51").Value = '  +2.62%  '
